$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156; everything from the old row 156 onward
# (through the old last row, 204) shifts down by one, becoming rows
# 157-205, and the workbook dimension grows from A1:T204 to A1:T205.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new price-report record.
$ws.Cells.Item(156, 1).Value = 2
$ws.Cells.Item(156, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(156, 3).Value = "Coquimbo"
$ws.Cells.Item(156, 4).Value = 44636
$ws.Cells.Item(156, 5).Value = 4
$ws.Cells.Item(156, 6).Value = "Fruta"
$ws.Cells.Item(156, 7).Value = 100103
$ws.Cells.Item(156, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(156, 9).Value = 100103006
$ws.Cells.Item(156, 10).Value = "Nectarín"
$ws.Cells.Item(156, 11).Value = "Artic Snow"
$ws.Cells.Item(156, 12).Value = "Especial"
$ws.Cells.Item(156, 13).Value = 16
$ws.Cells.Item(156, 14).Value = 390000
$ws.Cells.Item(156, 15).Value = 400000
$ws.Cells.Item(156, 16).Value = 395000
$ws.Cells.Item(156, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(156, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(156, 19).Value = 940
$ws.Cells.Item(156, 20).Value = 420
